$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting for the new column L (2021) from the corresponding column K (2020) cells,
# then set the new values for 2021.

# Row 3: empty bottom-border cell
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)

# Row 4: year header
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value = 2021

# Row 6: Mammals value
$ws.Range("K6").Copy()
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("L6").Value = 7.1

# Row 7: Birds value
$ws.Range("K7").Copy()
$ws.Range("L7").PasteSpecial(-4122)
$ws.Range("L7").Value = 0.5

# Row 8: Amphibians and Reptiles value -> "-" (same as the rest of the row)
$ws.Range("K8").Copy()
$ws.Range("L8").PasteSpecial(-4122)
$ws.Range("L8").Value = "-"

$excel.CutCopyMode = 0

$ws.Range("N5").Select()
